# Update Name of Algo
# Apply corrected values in column E (header "F") for the RandomForest result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.7897
$ws.Range("E10").Value = 12.30029999999999
$ws.Range("E12").Value = 12.147
$ws.Range("E18").Value = 12.1732
$ws.Range("E25").Value = 13.40419999999999
